$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells that hold numeric-looking text stay as TEXT
# (matches original inlineStr formatting, e.g. "1.00", "0.999") instead of
# being auto-converted to numbers by Excel.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.386.19"
$ws.Range("E2").Value = "  +4.36%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.496.45"
$ws.Range("E3").Value = "  +5.68%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.27%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.12"
$ws.Range("E5").Value = "  +7.54%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "184.92"
$ws.Range("E6").Value = "  +8.35%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.643"
$ws.Range("E7").Value = "  +9.43%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.489.58"
$ws.Range("E8").Value = "  +5.49%  "

# Row 9
$ws.Range("E9").Value = "  +0.02%  "

# Row 10
$ws.Range("E10").Value = "  +5.83%  "

# Row 11
$ws.Range("E11").Value = "  +16.44%  "

# Row 12
$ws.Range("E12").Value = "  +3.68%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000272"
$ws.Range("E13").Value = "  +7.27%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.28"
$ws.Range("E14").Value = "  +4.29%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.046.55"
$ws.Range("E15").Value = "  +4.64%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.492.25"
$ws.Range("E16").Value = "  +4.99%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.54"
$ws.Range("E17").Value = "  +6.98%  "

# Row 18
$ws.Range("E18").Value = "  +4.83%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "66.386.07"
$ws.Range("E19").Value = "  +4.54%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.05"
$ws.Range("E20").Value = "  +8.90%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.993"
$ws.Range("E21").Value = "  +4.87%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "422.03"
$ws.Range("E22").Value = "  +13.65%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.07"
$ws.Range("E23").Value = "  +11.91%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.20"
$ws.Range("E24").Value = "  +6.53%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.14"
$ws.Range("E25").Value = "  -0.90%  "

# Row 26
$ws.Range("E26").Value = "  -2.82%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.90"
$ws.Range("E27").Value = "  +9.02%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.33"
$ws.Range("E28").Value = "  +10.96%  "

# Row 29
$ws.Range("E29").Value = "  -1.39%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.13"
$ws.Range("E30").Value = "  +13.12%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.16"
$ws.Range("E31").Value = "  +5.89%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "630.56"
$ws.Range("E32").Value = "  +0.91%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.59"
$ws.Range("E33").Value = "  +4.36%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.72"
$ws.Range("E34").Value = "  +5.87%  "

# Row 35
$ws.Range("E35").Value = "  +6.44%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "59.87"
$ws.Range("E36").Value = "  +4.04%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0809"
$ws.Range("E37").Value = "  +12.73%  "

# Row 38 (reordered)
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.147"
$ws.Range("E38").Value = "  +19.83%  "

# Row 39 (reordered)
$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.03%  "

# Row 40 (reordered)
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.58"
$ws.Range("E40").Value = "  +5.88%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.385"
$ws.Range("E41").Value = "  +3.45%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.49"
$ws.Range("E42").Value = "  +16.94%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.117.93"
$ws.Range("E43").Value = "  +7.94%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("E44").Value = "  -0.40%  "

# Row 45
$ws.Range("E45").Value = "  -0.08%  "

# Row 46 (reordered)
$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.84"
$ws.Range("E46").Value = "  +11.07%  "

# Row 47 (reordered)
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.31"
$ws.Range("E47").Value = "  +12.12%  "

# Row 48
$ws.Range("E48").Value = "  +6.05%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.74"
$ws.Range("E49").Value = "  +2.61%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.134"
$ws.Range("E50").Value = "  +8.30%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "139.34"
$ws.Range("E51").Value = "  +2.88%  "
